$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "27.564.99"
Set-TextValue "D3" "1.724.22"
Set-TextValue "E3" "  +4.47%  "
Set-TextValue "D4" "1.004"
Set-TextValue "E4" "  +0.16%  "
Set-TextValue "D5" "226.09"
Set-TextValue "E5" "  +3.57%  "
Set-TextValue "D6" "0.5361"
Set-TextValue "E6" "  +3.10%  "
Set-TextValue "D7" "1.004"
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "0.2671"
Set-TextValue "E8" "  +1.31%  "
Set-TextValue "D9" "0.06607"
Set-TextValue "D10" "21.77"
Set-TextValue "E10" "  +6.76%  "
Set-TextValue "D11" "0.07741"
Set-TextValue "E11" "  +1.14%  "
Set-TextValue "D12" "4.626"
Set-TextValue "E12" "  +0.98%  "
Set-TextValue "D13" "1.722.53"
Set-TextValue "E13" "  +4.44%  "
Set-TextValue "D14" "1.962.42"
Set-TextValue "E14" "  +4.55%  "
Set-TextValue "D15" "0.5855"
Set-TextValue "E15" "  +4.46%  "
Set-TextValue "D16" "0.0₅8322"
Set-TextValue "E16" "  +2.04%  "
Set-TextValue "D17" "67.96"
Set-TextValue "E17" "  +4.22%  "
Set-TextValue "D18" "27.577.61"
Set-TextValue "E18" "  +5.90%  "
Set-TextValue "D19" "220.90"
Set-TextValue "E19" "  +15.26%  "
Set-TextValue "E20" "  +0.05%  "
Set-TextValue "D21" "4.730"
Set-TextValue "E21" "  +2.38%  "
Set-TextValue "D22" "10.64"
Set-TextValue "E22" "  +1.41%  "
Set-TextValue "D23" "6.094"
Set-TextValue "D24" "1.005"
Set-TextValue "E24" "  +0.12%  "
Set-TextValue "D25" "148.51"
Set-TextValue "E25" "  +3.25%  "
Set-TextValue "D26" "1.733"
Set-TextValue "E26" "  +14.58%  "
Set-TextValue "D27" "0.1237"
Set-TextValue "E27" "  +4.25%  "
Set-TextValue "D28" "7.421"
Set-TextValue "E28" "  +3.10%  "
Set-TextValue "D29" "16.60"
Set-TextValue "E29" "  +4.59%  "
Set-TextValue "D30" "0.05580"
Set-TextValue "E30" "  +3.19%  "
Set-TextValue "D31" "1.305"
Set-TextValue "E31" "  +2.91%  "
Set-TextValue "D32" "3.576"
Set-TextValue "E32" "  +3.76%  "
Set-TextValue "D33" "3.451"
Set-TextValue "E33" "  +3.01%  "
Set-TextValue "D34" "1.664"
Set-TextValue "E34" "  +7.09%  "
Set-TextValue "B35" "MXToken"
Set-TextValue "C35" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D35" "2.855"
Set-TextValue "E35" "  +2.63%  "
Set-TextValue "B36" "ARBITRUM"
Set-TextValue "C36" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D36" "0.9694"
Set-TextValue "E36" "  +2.56%  "
Set-TextValue "D37" "2.426"
Set-TextValue "E37" "  +0.04%  "
Set-TextValue "D38" "0.5975"
Set-TextValue "E38" "  +6.23%  "
Set-TextValue "D39" "0.01655"
Set-TextValue "E39" "  +4.91%  "
Set-TextValue "D40" "5.925"
Set-TextValue "E40" "  +1.00%  "
Set-TextValue "D41" "1.057.17"
Set-TextValue "E41" "  +2.84%  "
Set-TextValue "D42" "0.8539"
Set-TextValue "E42" "  +3.33%  "
Set-TextValue "E43" "  +0.10%  "
Set-TextValue "E44" "  +0.26%  "
Set-TextValue "D45" "1.868.39"
Set-TextValue "E45" "  +4.58%  "
Set-TextValue "D46" "0.0₈116"
Set-TextValue "E46" "  +3.98%  "
Set-TextValue "D47" "59.12"
Set-TextValue "E47" "  +2.97%  "
Set-TextValue "D48" "8.264"
Set-TextValue "E48" "  +4.47%  "
Set-TextValue "D49" "0.4436"
Set-TextValue "E49" "  +2.70%  "
Set-TextValue "E50" "  +0.74%  "
Set-TextValue "D51" "0.05255"
Set-TextValue "E51" "  +2.29%  "
